$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmails = @(
    "scam@mcafee.com  ",
    "neil_tyagi@mcafee.com  ",
    "somethingsomething@mcafee.com  ",
    "valerie_lenihan@mcafee.com",
    "stephanie.jones1@motorolasolutions.com  ",
    "ichrak.mekni@motorolasolutions.com  ",
    "shrinidhi.kannan@motorolasolutions.com  ",
    "john.c@motorolasolutions.com  ",
    "sayantan.karmakar@motorolasolutions.com  ",
    "ashish.jaiswal1@motorolasolutions.com  ",
    "tomer.dobershtein@motorolasolutions.com  ",
    "edward.connolly@motorolasolutions.com  ",
    "callum.crowe@motorolasolutions.com  ",
    "ecosystemproservices@motorolasolutions.com  ",
    "gille.smits@motorolasolutions.com  ",
    "dario.radosevic@motorolasolutions.com  ",
    "liviu.mandru@motorolasolutions.com  ",
    "silentsentinel.infoservice@motorolasolutions.com  ",
    "leon.dsouza@motorolasolutions.com  ",
    "justin.macdaniel@motorolasolutions.com  ",
    "anton.bouwer@motorolasolutions.com  ",
    "consultant@motorolasolutions.com  ",
    "marsha.tart@motorolasolutions.com  ",
    "ATInfo@motorolasolutions.com  ",
    "offthebeat@motorolasolutions.com  ",
    "mark.coughlan@motorolasolutions.com  ",
    "kruba.andalnesan@motorolasolutions.com  ",
    "drnaffiliate.managers@motorolasolutions.com  ",
    "malgorzata.duniec@motorolasolutions.com  ",
    "Jaroslaw.Magera@motorolasolutions.com  ",
    "jeff.corr@motorolasolutions.com",
    "careers@darwinbox.com  ",
    "sapna.rani@acuitykp.com  ",
    "nikhil.tyagi@acuitykp.com  ",
    "anandan.selvam@o9solutions.com  ",
    "C_shilpa.narayana@o9solutions.com  ",
    "reddy.babu@o9solutions.com  ",
    "c_suchismita.swami@o9solutions.com  ",
    "kausik.pradhan@o9solutions.com  ",
    "kusuma.manjunath@o9solutions.com  ",
    "vignesh.venkatesan@o9solutions.com  ",
    "minerva@o9solutions.com  ",
    "yuvraj.uppal@mapmyindia.com  ",
    "piyushmittal92555@gmail.com  "
)

$startRow = 3397
for ($i = 0; $i -lt $newEmails.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newEmails[$i]
}
